$d = $word.ActiveDocument

$d.Content.Find.Execute(
  "SISTEMA DE ARCHIVOS BASADO EN PLATAFORMA WEB PARA EVITAR LA FILTRACION DE DOCUMENTACION  CONFIDENCIAL  APLICANDO INTELIGENCIA  ARTIFICIAL, EXPRESIONES REGULARES Y ",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "SISTEMA DE ARCHIVOS BASADO  EN PLATAFORMA WEB PARA EVITAR LA FILTRACION  DE DOCUMENTACION CONFIDENCIAL APLICANDO INTELIGENCIA ARTIFICIAL,  EXPRESIONES REGULARES Y ",
  2)

$d.Content.Find.Execute(
  "Desarrollar un sistema de archivos de forma que ayude a la empresa  a evitar la filtraci",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Desarrollar  un sistema de archivos de forma que ayude a la empresa a evitar la filtraci",
  2)

$d.Content.Find.Execute(
  " otros. Haciendo uso de marcado esteganogr",
  $true, $false, $false, $false, $false, $true, 1, $false,
  " otros. Haciendo  uso de marcado esteganogr",
  2)

$d.Content.Find.Execute(
  "fico en los documentos  de forma que se pueda identificar al personal que realizo la filtraci",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "fico en los documentos de forma que se pueda identificar al personal que realizo la filtraci",
  2)

$d.Content.Find.Execute(
  "ginas de los documentos filtrados  puedan ser identificados por el mismo sistema, ya sea si estos son manipulados  en formato digital o en formato f",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "ginas de los documentos filtrados puedan ser identificados por el mismo sistema,  ya sea si estos son manipulados en formato digital o en formato f",
  2)

$d.Content.Find.Execute(
  "gina para que el sistema pueda identificar la procedencia del documento filtrado previamente marcado por el sistema. El marcado esteganogr",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "gina para que el sistema pueda identificar la procedencia  del documento filtrado previamente marcado por el sistema. El marcado esteganogr",
  2)

$d.Content.Find.Execute(
  " en el documento incluso si este es impreso  y manipulado de forma f",
  $true, $false, $false, $false, $false, $true, 1, $false,
  " en el documento incluso  si este es impreso y manipulado de forma f",
  2)

$d.Content.Find.Execute(
  "Para que el sistema pueda identificar la procedencia de cualquier documentaci",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Para que el sistema pueda  identificar la procedencia  de cualquier documentaci",
  2)

$d.Content.Find.Execute(
  "a en forma to JPG o el mismo documento en forma to PDF, el sistema evaluara la informa ci",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "a en formato JPG o el mismo documento en formato PDF, el sistema evaluara la informaci",
  2)

$d.Content.Find.Execute(
  "n obtenida  haciendo uso de inteligencia  artificial, expresiones regulares y esteganograf",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "n obtenida haciendo uso de inteligencia artificial,  expresiones regulares y esteganograf",
  2)

$d.Content.Find.Execute(
  "a. De esta forma  evitar la filtraci",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "a. De esta forma evitar la filtraci",
  2)
